$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A2").Value = 4845
$ws.Range("B2").Value = 7408
$ws.Range("C2").Value = 12253
$ws.Range("D2").Value = 2160.618445841338
$ws.Range("E2").Value = -964.7371519586635
$ws.Range("F2").Value = 0.2171133938593462
$ws.Range("G2").Value = -0.09550606636550885
$ws.Range("H2").Value = 0.02810784735271694
$ws.Range("I2").Value = 10468196.37010127
$ws.Range("J2").Value = -7146772.821709784
$ws.Range("K2").Value = 0.3954133681547376
$ws.Range("L2").Value = 2.239592868850058
$ws.Range("M2").Value = 1.464744526131009
$ws.Range("N2").Value = 3321423.548391489
